# AutoCommit_21 июня 2024 г. 9:22:01_SibNout2023
# Adds a "Вар" (variant) column in T, fills in variant numbers for some
# students, and moves the frozen-pane scroll / active selection to the
# newly edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Вар" column (row 2) - becomes a new shared string.
$ws.Range("T2").Value = "Вар"

# Variant numbers for a handful of students (rows 8,14,16,18,19,22,23,27).
$ws.Range("T8").Value = 18
$ws.Range("T14").Value = 23
$ws.Range("T16").Value = 25
$ws.Range("T18").Value = 3
$ws.Range("T19").Value = 11
$ws.Range("T22").Value = 9
$ws.Range("T23").Value = 16
$ws.Range("T27").Value = 12

# Move the active selection to the last touched cell.
$ws.Range("T25").Select()
